# Fixed a bug with precision/recall calculations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows after row 10 (pushes old row 12.. down to row 14..) ---
$ws.Rows("11:12").Insert()

# --- New column P width ---
$ws.Columns("P").ColumnWidth = 10.5

# --- Header cell for new column P ---
$ws.Range("P2").Value = "Pos#/Neg#"

# --- Fix experiment #7 (row 9): corrected color, sizes, and precision/recall figures ---
$ws.Range("D9").Value = "grayscale"
$ws.Range("E9").Value = 1900
$ws.Range("F9").Value = 1800
$ws.Range("L9").Value = 0.63
$ws.Range("M9").Value = 0.75
$ws.Range("N9").Value = 0.62
$ws.Range("O9").Value = "Cosine"
$ws.Range("P9").Value = 1.7

# --- Sheet view: selection moves to N9, scroll back to top ---
$ws.Range("N9").Select() | Out-Null
